$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$notes = $s.NotesPage
$notes.Shapes.Item(2).TextFrame.TextRange.Text = ""
